# Auto-generated edit script: re-run / update of RWheel run data (A:AD) and
# derived statistics (Fitness_SD, Fitness_Mean, Fitness_Lower, Fitness_Upper)
# for generations 0-10 (worksheet rows 2-12). Values correspond to a fresh
# simulation run ("temp solve of RWheel") pasted in as static data, matching
# the original authoring workflow (no formulas in the source file).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Build the 30-run matrix (columns A:AD) for rows 2-12 (generations 0-10)
$runData = New-Object 'object[,]' 11,30
$runData[0,0] = 3797
$runData[0,1] = 4187
$runData[0,2] = 3750
$runData[0,3] = 4098
$runData[0,4] = 3986
$runData[0,5] = 3274
$runData[0,6] = 4024
$runData[0,7] = 4274
$runData[0,8] = 4068
$runData[0,9] = 4092
$runData[0,10] = 3625
$runData[0,11] = 4329
$runData[0,12] = 4065
$runData[0,13] = 4091
$runData[0,14] = 4319
$runData[0,15] = 3917
$runData[0,16] = 4302
$runData[0,17] = 4058
$runData[0,18] = 4206
$runData[0,19] = 4357
$runData[0,20] = 4370
$runData[0,21] = 4171
$runData[0,22] = 3974
$runData[0,23] = 4004
$runData[0,24] = 4868
$runData[0,25] = 4308
$runData[0,26] = 4278
$runData[0,27] = 4297
$runData[0,28] = 3551
$runData[0,29] = 3960
$runData[1,0] = 3797
$runData[1,1] = 4187
$runData[1,2] = 3750
$runData[1,3] = 4098
$runData[1,4] = 3986
$runData[1,5] = 3872
$runData[1,6] = 4024
$runData[1,7] = 4297
$runData[1,8] = 4068
$runData[1,9] = 4092
$runData[1,10] = 3625
$runData[1,11] = 4471
$runData[1,12] = 4065
$runData[1,13] = 4091
$runData[1,14] = 4319
$runData[1,15] = 3917
$runData[1,16] = 4302
$runData[1,17] = 4058
$runData[1,18] = 4206
$runData[1,19] = 4357
$runData[1,20] = 4370
$runData[1,21] = 4171
$runData[1,22] = 3974
$runData[1,23] = 4004
$runData[1,24] = 4868
$runData[1,25] = 4308
$runData[1,26] = 4278
$runData[1,27] = 4297
$runData[1,28] = 3551
$runData[1,29] = 3960
$runData[2,0] = 3797
$runData[2,1] = 4187
$runData[2,2] = 3750
$runData[2,3] = 4098
$runData[2,4] = 3986
$runData[2,5] = 3872
$runData[2,6] = 4024
$runData[2,7] = 4297
$runData[2,8] = 4068
$runData[2,9] = 4092
$runData[2,10] = 3938
$runData[2,11] = 4471
$runData[2,12] = 4065
$runData[2,13] = 4091
$runData[2,14] = 4319
$runData[2,15] = 3917
$runData[2,16] = 4302
$runData[2,17] = 4058
$runData[2,18] = 4206
$runData[2,19] = 4357
$runData[2,20] = 4370
$runData[2,21] = 4171
$runData[2,22] = 3974
$runData[2,23] = 4217
$runData[2,24] = 4868
$runData[2,25] = 4308
$runData[2,26] = 4278
$runData[2,27] = 4297
$runData[2,28] = 3551
$runData[2,29] = 3960
$runData[3,0] = 3797
$runData[3,1] = 4187
$runData[3,2] = 4532
$runData[3,3] = 4098
$runData[3,4] = 3986
$runData[3,5] = 3872
$runData[3,6] = 4024
$runData[3,7] = 4297
$runData[3,8] = 4068
$runData[3,9] = 4092
$runData[3,10] = 3938
$runData[3,11] = 4471
$runData[3,12] = 4065
$runData[3,13] = 4091
$runData[3,14] = 4319
$runData[3,15] = 3917
$runData[3,16] = 4359
$runData[3,17] = 4058
$runData[3,18] = 4408
$runData[3,19] = 4357
$runData[3,20] = 4370
$runData[3,21] = 4171
$runData[3,22] = 3974
$runData[3,23] = 4217
$runData[3,24] = 4868
$runData[3,25] = 4308
$runData[3,26] = 4278
$runData[3,27] = 4297
$runData[3,28] = 3551
$runData[3,29] = 3960
$runData[4,0] = 3797
$runData[4,1] = 4187
$runData[4,2] = 4532
$runData[4,3] = 4098
$runData[4,4] = 3986
$runData[4,5] = 3872
$runData[4,6] = 4024
$runData[4,7] = 4297
$runData[4,8] = 4068
$runData[4,9] = 4092
$runData[4,10] = 3938
$runData[4,11] = 4471
$runData[4,12] = 4065
$runData[4,13] = 4091
$runData[4,14] = 4319
$runData[4,15] = 3917
$runData[4,16] = 4359
$runData[4,17] = 4058
$runData[4,18] = 4408
$runData[4,19] = 4357
$runData[4,20] = 4370
$runData[4,21] = 4171
$runData[4,22] = 3974
$runData[4,23] = 4217
$runData[4,24] = 4868
$runData[4,25] = 4467
$runData[4,26] = 4278
$runData[4,27] = 4297
$runData[4,28] = 3551
$runData[4,29] = 3960
$runData[5,0] = 3797
$runData[5,1] = 4187
$runData[5,2] = 4565
$runData[5,3] = 4098
$runData[5,4] = 3986
$runData[5,5] = 3872
$runData[5,6] = 4024
$runData[5,7] = 4297
$runData[5,8] = 4068
$runData[5,9] = 4416
$runData[5,10] = 3938
$runData[5,11] = 4471
$runData[5,12] = 4065
$runData[5,13] = 4091
$runData[5,14] = 4319
$runData[5,15] = 3917
$runData[5,16] = 4359
$runData[5,17] = 4058
$runData[5,18] = 4408
$runData[5,19] = 4357
$runData[5,20] = 4370
$runData[5,21] = 4171
$runData[5,22] = 3974
$runData[5,23] = 4217
$runData[5,24] = 4868
$runData[5,25] = 4467
$runData[5,26] = 4278
$runData[5,27] = 4297
$runData[5,28] = 3551
$runData[5,29] = 3960
$runData[6,0] = 3797
$runData[6,1] = 4187
$runData[6,2] = 4565
$runData[6,3] = 4098
$runData[6,4] = 3986
$runData[6,5] = 3872
$runData[6,6] = 4024
$runData[6,7] = 4297
$runData[6,8] = 4068
$runData[6,9] = 4416
$runData[6,10] = 3938
$runData[6,11] = 4471
$runData[6,12] = 4065
$runData[6,13] = 4091
$runData[6,14] = 4319
$runData[6,15] = 3917
$runData[6,16] = 4359
$runData[6,17] = 4058
$runData[6,18] = 4408
$runData[6,19] = 4357
$runData[6,20] = 4370
$runData[6,21] = 4171
$runData[6,22] = 3974
$runData[6,23] = 4217
$runData[6,24] = 4868
$runData[6,25] = 4467
$runData[6,26] = 4278
$runData[6,27] = 4297
$runData[6,28] = 3551
$runData[6,29] = 3960
$runData[7,0] = 3797
$runData[7,1] = 4187
$runData[7,2] = 4565
$runData[7,3] = 4098
$runData[7,4] = 4046
$runData[7,5] = 3872
$runData[7,6] = 4024
$runData[7,7] = 4297
$runData[7,8] = 4068
$runData[7,9] = 4416
$runData[7,10] = 3938
$runData[7,11] = 4471
$runData[7,12] = 4065
$runData[7,13] = 4091
$runData[7,14] = 4319
$runData[7,15] = 3917
$runData[7,16] = 4359
$runData[7,17] = 4058
$runData[7,18] = 4408
$runData[7,19] = 4357
$runData[7,20] = 4375
$runData[7,21] = 4171
$runData[7,22] = 3974
$runData[7,23] = 4217
$runData[7,24] = 4868
$runData[7,25] = 4467
$runData[7,26] = 4278
$runData[7,27] = 4297
$runData[7,28] = 3551
$runData[7,29] = 3960
$runData[8,0] = 3797
$runData[8,1] = 4187
$runData[8,2] = 4565
$runData[8,3] = 4098
$runData[8,4] = 4046
$runData[8,5] = 3872
$runData[8,6] = 4024
$runData[8,7] = 4297
$runData[8,8] = 4068
$runData[8,9] = 4416
$runData[8,10] = 3938
$runData[8,11] = 4471
$runData[8,12] = 4065
$runData[8,13] = 4091
$runData[8,14] = 4319
$runData[8,15] = 3917
$runData[8,16] = 4547
$runData[8,17] = 4058
$runData[8,18] = 4408
$runData[8,19] = 4357
$runData[8,20] = 4375
$runData[8,21] = 4171
$runData[8,22] = 3974
$runData[8,23] = 4217
$runData[8,24] = 4868
$runData[8,25] = 4467
$runData[8,26] = 4361
$runData[8,27] = 4297
$runData[8,28] = 3551
$runData[8,29] = 3960
$runData[9,0] = 3797
$runData[9,1] = 4187
$runData[9,2] = 4565
$runData[9,3] = 4098
$runData[9,4] = 4046
$runData[9,5] = 3872
$runData[9,6] = 4024
$runData[9,7] = 4297
$runData[9,8] = 4068
$runData[9,9] = 4416
$runData[9,10] = 3938
$runData[9,11] = 4471
$runData[9,12] = 4065
$runData[9,13] = 4091
$runData[9,14] = 4319
$runData[9,15] = 3917
$runData[9,16] = 4547
$runData[9,17] = 4058
$runData[9,18] = 4408
$runData[9,19] = 4357
$runData[9,20] = 4375
$runData[9,21] = 4171
$runData[9,22] = 3974
$runData[9,23] = 4217
$runData[9,24] = 4868
$runData[9,25] = 4663
$runData[9,26] = 4361
$runData[9,27] = 4297
$runData[9,28] = 3551
$runData[9,29] = 3960
$runData[10,0] = 3797
$runData[10,1] = 4187
$runData[10,2] = 4565
$runData[10,3] = 4098
$runData[10,4] = 4046
$runData[10,5] = 3872
$runData[10,6] = 4024
$runData[10,7] = 4297
$runData[10,8] = 4068
$runData[10,9] = 4416
$runData[10,10] = 3938
$runData[10,11] = 4471
$runData[10,12] = 4065
$runData[10,13] = 4091
$runData[10,14] = 4319
$runData[10,15] = 3917
$runData[10,16] = 4547
$runData[10,17] = 4058
$runData[10,18] = 4408
$runData[10,19] = 4357
$runData[10,20] = 4375
$runData[10,21] = 4171
$runData[10,22] = 3974
$runData[10,23] = 4217
$runData[10,24] = 4868
$runData[10,25] = 4663
$runData[10,26] = 4361
$runData[10,27] = 4297
$runData[10,28] = 3551
$runData[10,29] = 3960

$ws.Range("A2:AD12").Value = $runData

# Derived statistics: Fitness_SD (AF), Fitness_Mean (AG), Fitness_Lower (AH), Fitness_Upper (AI)
$statsData = New-Object 'object[,]' 11,4
$statsData[0,0] = 296.9529549798493
$statsData[0,1] = 4086.666666666667
$statsData[0,2] = 4383.619621646516
$statsData[0,3] = 3789.713711686817
$statsData[1,0] = 264.0293021878313
$statsData[1,1] = 4112.1
$statsData[1,2] = 4376.129302187832
$statsData[1,3] = 3848.070697812169
$statsData[2,0] = 249.466531973154
$statsData[2,1] = 4129.633333333333
$statsData[2,2] = 4379.099865306487
$statsData[2,3] = 3880.166801360179
$statsData[3,0] = 254.6791535724269
$statsData[3,1] = 4164.333333333333
$statsData[3,2] = 4419.01248690576
$statsData[3,3] = 3909.654179760906
$statsData[4,0] = 259.383019056207
$statsData[4,1] = 4169.633333333333
$statsData[4,2] = 4429.01635238954
$statsData[4,3] = 3910.250314277126
$statsData[5,0] = 264.3491334727922
$statsData[5,1] = 4181.533333333334
$statsData[5,2] = 4445.882466806126
$statsData[5,3] = 3917.184199860541
$statsData[6,0] = 264.3491334727922
$statsData[6,1] = 4181.533333333334
$statsData[6,2] = 4445.882466806126
$statsData[6,3] = 3917.184199860541
$statsData[7,0] = 263.1662841312471
$statsData[7,1] = 4183.7
$statsData[7,2] = 4446.866284131247
$statsData[7,3] = 3920.533715868753
$statsData[8,0] = 271.0000551384768
$statsData[8,1] = 4192.733333333334
$statsData[8,2] = 4463.73338847181
$statsData[8,3] = 3921.733278194857
$statsData[9,0] = 280.0515880226766
$statsData[9,1] = 4199.266666666666
$statsData[9,2] = 4479.318254689343
$statsData[9,3] = 3919.21507864399
$statsData[10,0] = 280.0515880226766
$statsData[10,1] = 4199.266666666666
$statsData[10,2] = 4479.318254689343
$statsData[10,3] = 3919.21507864399

$ws.Range("AF2:AI12").Value = $statsData

Write-Host "RWheel run data + fitness stats updated for rows 2-12"